$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5262
$ws1.Range("F6").Value = 5262
$ws1.Range("F7").Value = 155
$ws1.Range("F11").Value = 1184
$ws1.Range("F12").Value = 747
$ws1.Range("C13").Value = "北京·Aw动漫游戏嘉年华8th-夏日奇幻之旅（Part1）"
$ws1.Range("F13").Value = 5206
$ws1.Range("F17").Value = 2291
$ws1.Range("F18").Value = 2291
$ws1.Range("F19").Value = 257
$ws1.Range("F22").Value = 3914
$ws1.Range("F26").Value = 3840
$ws1.Range("F27").Value = 187
$ws1.Range("F36").Value = 22
$ws1.Range("G36").Value = 520
$ws1.Range("F37").Value = 6826
$ws1.Range("F38").Value = 1107
$ws1.Range("F39").Value = 528
$ws1.Range("F42").Value = 1388
$ws1.Range("F44").Value = 696
$ws1.Range("F46").Value = 2321
$ws1.Range("F50").Value = 786
$ws1.Range("F51").Value = 931

# Sheet "演出" (index 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 8
$ws2.Range("F22").Value = 53
$ws2.Range("F25").Value = 819

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 5262
$ws4.Range("F8").Value = 5262
$ws4.Range("F9").Value = 155
$ws4.Range("F14").Value = 1184
$ws4.Range("F15").Value = 747
$ws4.Range("F19").Value = 2291
$ws4.Range("F20").Value = 2291
$ws4.Range("F21").Value = 257
$ws4.Range("F24").Value = 3914
$ws4.Range("F25").Value = 3840
$ws4.Range("F26").Value = 187
$ws4.Range("F34").Value = 22
$ws4.Range("G34").Value = 520
$ws4.Range("F36").Value = 6826
$ws4.Range("F37").Value = 1107
$ws4.Range("F38").Value = 528
$ws4.Range("F42").Value = 1388
$ws4.Range("F44").Value = 696
$ws4.Range("F46").Value = 2322
$ws4.Range("F49").Value = 786
$ws4.Range("F50").Value = 931
